$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E, shifting existing E:I -> F:J
$ws.Columns("E").Insert()

# New column header
$ws.Range("E1").Value = "domestic_final_demand_index"

# Rename the (shifted) 3yr bond spread header to 3yr credit spread
$ws.Range("G1").Value = "3yr_credit_spread_qrtly"

# Populate the new domestic_final_demand_index column for rows 2-65
$dfdiValues = @(69, 69.3, 70, 70.5, 71.09999999999999, 71.90000000000001, 72.59999999999999, 73, 73.3, 73.90000000000001, 74.59999999999999, 75.09999999999999, 75.90000000000001, 76.59999999999999, 77.8, 78.7, 79.2, 79.3, 79.7, 79.90000000000001, 80.40000000000001, 80.7, 81.40000000000001, 81.59999999999999, 82.2, 82.5, 83, 83.3, 83.40000000000001, 83.90000000000001, 84.7, 84.90000000000001, 85.09999999999999, 85.40000000000001, 86.40000000000001, 86.90000000000001, 87.40000000000001, 87.59999999999999, 88.09999999999999, 88.5, 88.90000000000001, 89.3, 89.90000000000001, 90.2, 90, 90, 90.3, 90.59999999999999, 90.8, 91.09999999999999, 91.40000000000001, 91.90000000000001, 92.3, 92.7, 93.3, 93.8, 94, 94.5, 94.90000000000001, 95.40000000000001, 95.59999999999999, 95.5, 95.59999999999999, 96)
$row = 2
foreach ($v in $dfdiValues) {
    $ws.Cells.Item($row, 5).Value = $v
    $row = $row + 1
}
